$wb = $excel.ActiveWorkbook

# Status text that changed from "Ready for handoff" to "Handed back: in sync with en-US"
# appears on both locale sheets, rows 2 and 3 (column B).
$newStatus = "Handed back: in sync with en-US"

$localeSheets = @(
  @{ Sheet = $wb.Worksheets.Item("zh-cn"); HandbackTime = "2016-03-11 02:12:35" },
  @{ Sheet = $wb.Worksheets.Item("de-de"); HandbackTime = "2016-03-11 02:12:58" }
)

foreach ($entry in $localeSheets) {
  $ws = $entry.Sheet
  $handbackTime = $entry.HandbackTime

  for ($row = 2; $row -le 3; $row++) {
    # Column B: Status -> report handed back and in sync with en-US source.
    $ws.Cells.Item($row, 2).Value = $newStatus

    # Column A holds the source file (.md) hyperlink/display text; Column C
    # holds the latest handoff (.xlf) hyperlink/display text. The handback
    # step fills in "Latest Target File" (E) and "Latest Handback File" (F)
    # with the very same files, each carrying its own hyperlink.
    $sourceCell = $ws.Cells.Item($row, 1)
    $handoffCell = $ws.Cells.Item($row, 3)

    $sourceText = $sourceCell.Value()
    $handoffText = $handoffCell.Value()

    $sourceLink = $sourceCell.Hyperlinks.Item(1)
    $handoffLink = $handoffCell.Hyperlinks.Item(1)

    $targetCell = $ws.Cells.Item($row, 5)
    $handbackCell = $ws.Cells.Item($row, 6)

    $ws.Hyperlinks.Add($targetCell, $sourceLink.Address, "", "", $sourceText)
    $ws.Hyperlinks.Add($handbackCell, $handoffLink.Address, "", "", $handoffText)

    # Column G: Latest Handback DateTime -> now populated with the real
    # handback timestamp instead of the 0001-01-01 placeholder.
    $ws.Cells.Item($row, 7).Value = $handbackTime
  }
}
